$d = $word.ActiveDocument

# Title change (appears twice: Heading1 at top, and bold text near bottom)
$d.Content.Find.Execute("Play Esqueleto Mariachi Slot for Free | Review", $true, $false, $false, $false, $false, $true, 1, $false, "Play Esqueleto Mariachi for Free - Exciting Slot Game Review", 2)
$d.Content.Find.Execute("Play Esqueleto Mariachi Slot for Free | Review", $true, $false, $false, $false, $false, $true, 1, $false, "Play Esqueleto Mariachi for Free - Exciting Slot Game Review", 2)

# "What we like" bullet list
$d.Content.Find.Execute("Captivating graphics and sound", $true, $false, $false, $false, $false, $true, 1, $false, "Captivating graphics and well-realized scenes", 2)
$d.Content.Find.Execute("Fun and engaging theme with unique symbols", $true, $false, $false, $false, $false, $true, 1, $false, "Engaging theme with vibrant colors and characteristic images", 2)
$d.Content.Find.Execute("Special features and free spins increase chances of winning", $true, $false, $false, $false, $false, $true, 1, $false, "Exciting special features and bonuses", 2)
$d.Content.Find.Execute("Maximum payout of €500,000", $true, $false, $false, $false, $false, $true, 1, $false, "High maximum payout potential", 2)

# "What we don't like" bullet list
$d.Content.Find.Execute("Fixed paylines limit betting options", $true, $false, $false, $false, $false, $true, 1, $false, "Limited betting range for high rollers", 2)
$d.Content.Find.Execute("Limited bonus games", $true, $false, $false, $false, $false, $true, 1, $false, "Limited number of paylines", 2)

# Meta description (italic text)
$d.Content.Find.Execute("Read our review of Esqueleto Mariachi slot game. Play this online casino game for free and win big prizes with its engaging theme and unique special features.", $true, $false, $false, $false, $false, $true, 1, $false, "Read our review of Esqueleto Mariachi, a fun and engaging slot game. Play for free and enjoy exciting special features.", 2)
